$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.754972666666667
$ws.Range("H2").Value = 5.264918
$ws.Range("I2").Value = 0.5110994274238188
$ws.Range("J2").Value = 0.5110994274238188
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 29.01953304682399
$ws.Range("R2").Value = 261.175797421416
$ws.Range("S2").Value = 0.1083762389050867
$ws.Range("T2").Value = 0.1083762389050867

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.754972666666667
$ws.Range("H3").Value = 5.264918
$ws.Range("I3").Value = 0.5110994274238188
$ws.Range("J3").Value = 0.5110994274238188
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 71.28809652276178
$ws.Range("R3").Value = 641.5928687048561
$ws.Range("S3").Value = 0.266232256989582
$ws.Range("T3").Value = 0.266232256989582

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.754972666666667
$ws.Range("H4").Value = 5.264918
$ws.Range("I4").Value = 0.5110994274238188
$ws.Range("J4").Value = 0.5110994274238188
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 36.547707672074
$ws.Range("R4").Value = 328.929369048666
$ws.Range("S4").Value = 0.1364909315291501
$ws.Range("T4").Value = 0.1364909315291501

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.678748
$ws.Range("H5").Value = 5.036244
$ws.Range("I5").Value = 0.4889005725761812
$ws.Range("J5").Value = 0.4889005725761812
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 27.759112143792
$ws.Range("R5").Value = 249.832009294128
$ws.Range("S5").Value = 0.1036690757440685
$ws.Range("T5").Value = 0.1036690757440685

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.678748
$ws.Range("H6").Value = 5.036244
$ws.Range("I6").Value = 0.4889005725761812
$ws.Range("J6").Value = 0.4889005725761812
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 68.19180249040534
$ws.Range("R6").Value = 613.7262224136481
$ws.Range("S6").Value = 0.2546688489488802
$ws.Range("T6").Value = 0.2546688489488802

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.678748
$ws.Range("H7").Value = 5.036244
$ws.Range("I7").Value = 0.4889005725761812
$ws.Range("J7").Value = 0.4889005725761812
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 34.960311533292
$ws.Range("R7").Value = 314.642803799628
$ws.Range("S7").Value = 0.1305626478832325
$ws.Range("T7").Value = 0.1305626478832326
